# edit.ps1 — apply sumN.xlsx update
# 1) Swap the "sumN_opt2" / "sumN_recursion" column headers (C1 <-> D1)
# 2) Replace the benchmark timing data in B2:H18 with the new run's numbers
# 3) Re-size columns B:H to fit the new header text
# 4) Leave the selection on E25, matching the saved workbook view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap header text for columns C and D ---
$ws.Range("C1").Value = "sumN_recursion"
$ws.Range("D1").Value = "sumN_opt2"

# --- Column widths (characters) for B:H ---
$ws.Columns.Item(2).ColumnWidth = 13.0
$ws.Columns.Item(3).ColumnWidth = 14.0
$ws.Columns.Item(4).ColumnWidth = 19.714285714285715
$ws.Columns.Item(5).ColumnWidth = 19.857142857142858
$ws.Columns.Item(6).ColumnWidth = 25.0
$ws.Columns.Item(7).ColumnWidth = 15.0
$ws.Columns.Item(8).ColumnWidth = 11.428571428571429

# --- Updated benchmark data (rows 2-18, columns B-H) ---
# Row 2
$ws.Range("B2").Value = 0.028367
$ws.Range("C2").Value = 0.031309999999999998
$ws.Range("D2").Value = 0.030498000000000001
$ws.Range("E2").Value = 0.044838000000000003
$ws.Range("F2").Value = 0.050243000000000003
$ws.Range("G2").Value = 0.047885999999999998
$ws.Range("H2").Value = 0.050776000000000002
# Row 3
$ws.Range("B3").Value = 0.053261999999999997
$ws.Range("C3").Value = 0.057504
$ws.Range("D3").Value = 0.054300000000000001
$ws.Range("E3").Value = 0.081394999999999995
$ws.Range("F3").Value = 0.088100999999999999
$ws.Range("G3").Value = 0.085565000000000002
$ws.Range("H3").Value = 0.094946000000000003
# Row 4
$ws.Range("B4").Value = 0.119212
$ws.Range("C4").Value = 0.11506
$ws.Range("D4").Value = 0.101204
$ws.Range("E4").Value = 0.13893
$ws.Range("F4").Value = 0.14475199999999999
$ws.Range("G4").Value = 0.15259400000000001
$ws.Range("H4").Value = 0.143675
# Row 5
$ws.Range("B5").Value = 0.221439
$ws.Range("C5").Value = 0.25043199999999999
$ws.Range("D5").Value = 0.16255900000000001
$ws.Range("E5").Value = 0.24957699999999999
$ws.Range("F5").Value = 0.26889800000000003
$ws.Range("G5").Value = 0.28089900000000001
$ws.Range("H5").Value = 0.27535500000000002
# Row 6
$ws.Range("B6").Value = 0.45430999999999999
$ws.Range("C6").Value = 0.44952599999999998
$ws.Range("D6").Value = 0.30221300000000001
$ws.Range("E6").Value = 0.46366800000000002
$ws.Range("F6").Value = 0.42622599999999999
$ws.Range("G6").Value = 0.50500400000000001
$ws.Range("H6").Value = 0.45693400000000001
# Row 7
$ws.Range("B7").Value = 0.90490099999999996
$ws.Range("C7").Value = 0.93520099999999995
$ws.Range("D7").Value = 0.54765799999999998
$ws.Range("E7").Value = 0.89277099999999998
$ws.Range("F7").Value = 1.1113900000000001
$ws.Range("G7").Value = 0.77101799999999998
$ws.Range("H7").Value = 0.79245600000000005
# Row 8
$ws.Range("B8").Value = 1.8161700000000001
$ws.Range("C8").Value = 2.0171299999999999
$ws.Range("D8").Value = 0.97164499999999998
$ws.Range("E8").Value = 1.8588899999999999
$ws.Range("F8").Value = 2.3436499999999998
$ws.Range("G8").Value = 1.40944
$ws.Range("H8").Value = 1.4699199999999999
# Row 9
$ws.Range("B9").Value = 4.1983300000000003
$ws.Range("C9").Value = 3.52841
$ws.Range("D9").Value = 1.9854499999999999
$ws.Range("E9").Value = 3.57613
$ws.Range("F9").Value = 4.1147200000000002
$ws.Range("G9").Value = 2.7212900000000002
$ws.Range("H9").Value = 2.86069
# Row 10
$ws.Range("B10").Value = 8.3974100000000007
$ws.Range("C10").Value = 6.3203899999999997
$ws.Range("D10").Value = 2.6357900000000001
$ws.Range("E10").Value = 7.1013700000000002
$ws.Range("F10").Value = 7.4751599999999998
$ws.Range("G10").Value = 5.3988100000000001
$ws.Range("H10").Value = 5.7024400000000002
# Row 11
$ws.Range("B11").Value = 17.368600000000001
$ws.Range("C11").Value = 11.9194
$ws.Range("D11").Value = 5.2257999999999996
$ws.Range("E11").Value = 14.075900000000001
$ws.Range("F11").Value = 13.841799999999999
$ws.Range("G11").Value = 10.5878
$ws.Range("H11").Value = 11.0337
# Row 12
$ws.Range("B12").Value = 33.866199999999999
$ws.Range("C12").Value = 23.473199999999999
$ws.Range("D12").Value = 10.236700000000001
$ws.Range("E12").Value = 28.150400000000001
$ws.Range("F12").Value = 27.076899999999998
$ws.Range("G12").Value = 20.874400000000001
$ws.Range("H12").Value = 21.034199999999998
# Row 13
$ws.Range("B13").Value = 67.133499999999998
$ws.Range("C13").Value = 43.9407
$ws.Range("D13").Value = 20.588999999999999
$ws.Range("E13").Value = 56.020299999999999
$ws.Range("F13").Value = 48.658999999999999
$ws.Range("G13").Value = 41.480600000000003
$ws.Range("H13").Value = 41.937800000000003
# Row 14
$ws.Range("B14").Value = 132.56200000000001
$ws.Range("C14").Value = 85.057400000000001
$ws.Range("D14").Value = 44.408099999999997
$ws.Range("E14").Value = 112.21
$ws.Range("F14").Value = 95.550200000000004
$ws.Range("G14").Value = 82.939099999999996
$ws.Range("H14").Value = 83.444000000000003
# Row 15
$ws.Range("B15").Value = 265.50700000000001
$ws.Range("C15").Value = 174.17500000000001
$ws.Range("D15").Value = 88.8904
$ws.Range("E15").Value = 225.066
$ws.Range("F15").Value = 192.26400000000001
$ws.Range("G15").Value = 167.64
$ws.Range("H15").Value = 170.18600000000001
# Row 16
$ws.Range("B16").Value = 533.45699999999999
$ws.Range("C16").Value = 345.95400000000001
$ws.Range("D16").Value = 165.512
$ws.Range("E16").Value = 451.666
$ws.Range("F16").Value = 382.69299999999998
$ws.Range("G16").Value = 345.79
$ws.Range("H16").Value = 341.64499999999998
# Row 17
$ws.Range("B17").Value = 1075.8
$ws.Range("C17").Value = 657.76099999999997
$ws.Range("D17").Value = 335.45699999999999
$ws.Range("E17").Value = 902.78099999999995
$ws.Range("F17").Value = 760.51599999999996
$ws.Range("G17").Value = 673.00199999999995
$ws.Range("H17").Value = 687.178
# Row 18
$ws.Range("B18").Value = 2058.7199999999998
$ws.Range("C18").Value = 1284.25
$ws.Range("D18").Value = 675.822
$ws.Range("E18").Value = 1846.27
$ws.Range("F18").Value = 1514.17
$ws.Range("G18").Value = 1364.04
$ws.Range("H18").Value = 1516.47

# --- Update the active selection to match the saved view ---
$ws.Range("E25").Select() | Out-Null
